# Regenerate the "K" (strikes) column (column G) values on the save_data
# sheet. The source data pipeline recomputed these values (switching from
# a raw "Strike#" count to the new "K" metric), so here we just write the
# newly calculated values into the existing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 3
    6  = 6
    7  = 8
    8  = 1
    9  = 5
    10 = 6
    11 = 4
    12 = 2
    13 = 6
    14 = 6
    15 = 1
    16 = 3
    17 = 1
    18 = 2
    19 = 0
    20 = 4
    21 = 2
    22 = 6
    23 = 0
    24 = 0
    25 = 5
    28 = 1
    29 = 3
    30 = 1
    31 = 4
    32 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
